$d = $word.ActiveDocument

# Each entry describes a table cell whose text like "26:30 h" must become
# "29:30 h" etc. Only a single digit in the middle of the string actually
# changes; the original run gets split into three runs (prefix, changed
# digit, suffix) that all keep identical run formatting, matching how Word
# would represent an in-place single-character edit.
$edits = @(
    @{ Old = "26:30 h";  SplitAt = 1; NewChar = "9" },
    @{ Old = "55:50 h";  SplitAt = 1; NewChar = "8" },
    @{ Old = "83:40 h";  SplitAt = 1; NewChar = "6" },
    @{ Old = "173:15 h"; SplitAt = 2; NewChar = "6" }
)

foreach ($edit in $edits) {
    $old = $edit.Old
    $splitAt = $edit.SplitAt
    $newChar = $edit.NewChar

    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $start = $rng.Start
    $end = $rng.End

    # Replace the single middle character with its new value first.
    $midRange = $d.Range($start + $splitAt, $start + $splitAt + 1)
    $midRange.Text = $newChar

    # Now force the run to split into three pieces (prefix / changed char /
    # suffix) by touching (toggling) a formatting property on each piece.
    # Toggling Bold on then back off leaves the run boundaries in place
    # without altering the visible formatting.
    $r1 = $d.Range($start, $start + $splitAt)
    $r2 = $d.Range($start + $splitAt, $start + $splitAt + 1)
    $r3 = $d.Range($start + $splitAt + 1, $end)

    $r1.Font.Bold = 1
    $r1.Font.Bold = 0
    $r2.Font.Bold = 1
    $r2.Font.Bold = 0
    $r3.Font.Bold = 1
    $r3.Font.Bold = 0
}
